$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two inserted columns (T: gateway_port, U: http_port)
$ws.Range("T3").Value = "gateway_port"
$ws.Range("U3").Value = "http_port"

# Per-zone gateway_port / http_port values for rows 4-13
$ws.Range("T4").Value  = "20001"
$ws.Range("U4").Value  = "8888"

$ws.Range("T5").Value  = "20002"
$ws.Range("U5").Value  = "8889"

$ws.Range("T6").Value  = "20003"
$ws.Range("U6").Value  = "8890"

$ws.Range("T7").Value  = "20004"
$ws.Range("U7").Value  = "8891"

$ws.Range("T8").Value  = "20005"
$ws.Range("U8").Value  = "8892"

$ws.Range("T9").Value  = "20006"
$ws.Range("U9").Value  = "8893"

$ws.Range("T10").Value = "20007"
$ws.Range("U10").Value = "8894"

$ws.Range("T11").Value = "20008"
$ws.Range("U11").Value = "8895"

$ws.Range("T12").Value = "20009"
$ws.Range("U12").Value = "8896"

$ws.Range("T13").Value = "20010"
$ws.Range("U13").Value = "8897"

# Match the selection state recorded in the edit (active cell U4, selection U4:U13)
$ws.Range("U4:U13").Select()
